$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 425, shifting existing rows 425:548 down to 426:549.
$ws.Rows.Item(425).Insert()

# Populate the newly inserted row 425 with the new weekly price record.
$ws.Range("A425").Value() = 4
$ws.Range("B425").Value() = "Feria Lagunitas de Puerto Montt"
$ws.Range("C425").Value() = "Los Lagos"
$ws.Range("D425").Value() = 45093
$ws.Range("E425").Value() = 10
$ws.Range("F425").Value() = 100112023
$ws.Range("G425").Value() = "Brócoli"
$ws.Range("H425").Value() = "Sin especificar"
$ws.Range("I425").Value() = "Primera"
$ws.Range("J425").Value() = 1400
$ws.Range("K425").Value() = 1400
$ws.Range("L425").Value() = 1500
$ws.Range("M425").Value() = 1450
$ws.Range("N425").Value() = "$/unidad"
$ws.Range("O425").Value() = "Región Metropolitana"
$ws.Range("P425").Value() = 1450
$ws.Range("Q425").Value() = 1
$ws.Range("R425").Value() = "Hortaliza"
